$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 email domain
$ws.Range("C2").Value = "@gmail.com"

# Clear out rows 3 through 7 (these order rows are no longer present) and
# force the now-empty cells to stay materialized (matching the blank cell
# pattern already used by the untouched rows below them) instead of being
# dropped entirely from the sheet on save.
$rng = $ws.Range("A3:H7")
$rng.Value = $null
$rng.VerticalAlignment = -4107
